$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- 1 & 2: the long "player controls" paragraph (right after the bold
# "玩家：" heading) is shortened to an intro line, and a brand new paragraph
# carrying the detailed operations list is inserted right after it. ---
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "玩家可以使用方向键*") {
        $target = $para
        break
    }
}
if ($target -eq $null) {
    throw "Could not find the player-controls paragraph"
}

$introText = "玩家可以进行如下操作："
$opsText = "向左移动，向右移动，跳跃，下蹲，奔跑。部分关卡会限制玩家的移动速度以及奔跑功能"

$newParasXml = "<w:p $wNs>" +
    "<w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>$introText</w:t></w:r>" +
    "</w:p>" +
    "<w:p $wNs>" +
    "<w:pPr><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr></w:pPr>" +
    "<w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>$opsText</w:t></w:r>" +
    "</w:p>"

$target.Range.InsertXML($newParasXml)

# --- 3: the trailing empty paragraph (after "具体数值如下：") drops its
# stray paragraph-mark formatting, becoming a fully bare paragraph. ---
$last = $d.Paragraphs.Item($d.Paragraphs.Count)
$last.Range.InsertXML("<w:p $wNs/>")

Write-Output "ok"
